# TC05_CDS_Filter_PHSAccession-phs001787.xlsx - PHS Accession Filter Test cases - 13
# Update the Neo4j queries (StatQuery + per-tab queries) to the revised Cypher.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$statQuery = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001787"]
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001787"]
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001787"]
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

$filesQuery = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001787"]
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name limit 100
"@

$samplesQuery = @"
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001787"]
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(samp.sample_tumor_status,'') as ``Tumor``,
    coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER BY samp.sample_id limit 100
"@

$participantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001787"] 
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id limit 100
"@

# Row 2 = ParticipantsTab, Row 3 = SamplesTab, Row 4 = FilesTab
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# Row heights grow to the Excel max (409.5) because the new query text is
# much longer; row 2/3 previously autofit to 171.6, row 4 to 218.4.
$ws.Range("A2:E2").RowHeight = 409.5
$ws.Range("A3:E3").RowHeight = 409.5
$ws.Range("A4:E4").RowHeight = 409.5

# Selection moved from D2 to C2 (author was reviewing the StatQuery cell).
$ws.Range("C2").Select()
